$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E4").Value = 10.6
$ws.Range("F4").Value = 9.92

$ws.Range("D5").Value = 9.4
$ws.Range("F5").Value = 10.18

$ws.Range("D6").Value = 10.08
$ws.Range("E6").Value = 9.82
$ws.Range("G6").Value = 10.38

$ws.Range("F7").Value = 9.62
$ws.Range("H7").Value = 10.11

$ws.Range("G8").Value = 9.89
$ws.Range("J8").Value = 11.55

$ws.Range("H10").Value = 8.45
